$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new price values are plain numeric literals (single decimal
# point) must be kept as text to match the source data (the "Price" column
# is authored as text throughout, e.g. thousands-dotted values like
# "29.322.52" alongside plain decimals like "0.9985"). Force text format on
# those specific cells before writing so Excel does not reinterpret them as
# numbers (which would also silently drop trailing zeros, e.g. "0.6740").
$textCells = @(
    "D4",
    "D5",
    "D6",
    "D8",
    "D10",
    "D13",
    "D14",
    "D16",
    "D17",
    "D19",
    "D20",
    "D21",
    "D22",
    "D23",
    "D24",
    "D25",
    "D26",
    "D28",
    "D29",
    "D30",
    "D31",
    "D33",
    "D35",
    "D36",
    "D37",
    "D38",
    "D39",
    "D41",
    "D44",
    "D45",
    "D47",
    "D48",
    "D49",
    "D50"
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "29.322.52"
$ws.Range("E2").Value = "  -0.12%  "
$ws.Range("D3").Value = "1.838.83"
$ws.Range("E3").Value = "  -0.40%  "
$ws.Range("D4").Value = "0.9985"
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "239.27"
$ws.Range("E5").Value = "  -0.42%  "
$ws.Range("D6").Value = "0.6261"
$ws.Range("E6").Value = "  -0.18%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "0.07429"
$ws.Range("E8").Value = "  -0.76%  "
$ws.Range("E9").Value = "  -0.09%  "
$ws.Range("D10").Value = "24.82"
$ws.Range("E10").Value = "  +1.86%  "
$ws.Range("E11").Value = "  -0.21%  "
$ws.Range("D12").Value = "1.832.68"
$ws.Range("E12").Value = "  -0.71%  "
$ws.Range("D13").Value = "4.951"
$ws.Range("E13").Value = "  -0.90%  "
$ws.Range("D14").Value = "0.6740"
$ws.Range("E14").Value = "  -0.74%  "
$ws.Range("E15").Value = "  -1.02%  "
$ws.Range("D16").Value = "81.77"
$ws.Range("E16").Value = "  -0.39%  "
$ws.Range("D17").Value = "6.234"
$ws.Range("E17").Value = "  +1.16%  "
$ws.Range("D18").Value = "29.331.48"
$ws.Range("E18").Value = "  -0.26%  "
$ws.Range("D19").Value = "233.02"
$ws.Range("E19").Value = "  +1.74%  "
$ws.Range("D20").Value = "12.28"
$ws.Range("E20").Value = "  -0.32%  "
$ws.Range("D21").Value = "0.9999"
$ws.Range("E21").Value = "  +0.01%  "
$ws.Range("D22").Value = "7.343"
$ws.Range("E22").Value = "  -1.56%  "
$ws.Range("D23").Value = "0.9999"
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("D24").Value = "158.07"
$ws.Range("E24").Value = "  -0.33%  "
$ws.Range("D25").Value = "8.476"
$ws.Range("E25").Value = "  +0.83%  "
$ws.Range("D26").Value = "0.1346"
$ws.Range("E26").Value = "  -2.16%  "
$ws.Range("E27").Value = "  -1.01%  "
$ws.Range("D28").Value = "0.07318"
$ws.Range("E28").Value = "  +13.67%  "
$ws.Range("D29").Value = "1.455"
$ws.Range("E29").Value = "  +4.98%  "
$ws.Range("D30").Value = "1.477"
$ws.Range("E30").Value = "  +0.36%  "
$ws.Range("D31").Value = "4.038"
$ws.Range("E31").Value = "  -1.24%  "
$ws.Range("E32").Value = "  -0.57%  "
$ws.Range("D33").Value = "1.816"
$ws.Range("E33").Value = "  -0.59%  "
$ws.Range("E34").Value = "  -0.09%  "
$ws.Range("D35").Value = "0.6965"
$ws.Range("E35").Value = "  -0.56%  "
$ws.Range("D36").Value = "2.570"
$ws.Range("E36").Value = "  -0.35%  "
$ws.Range("D37").Value = "0.01832"
$ws.Range("E37").Value = "  +0.21%  "
$ws.Range("D38").Value = "2.815"
$ws.Range("E38").Value = "  -0.54%  "
$ws.Range("D39").Value = "6.870"
$ws.Range("E39").Value = "  +4.14%  "
$ws.Range("D40").Value = "1.230.44"
$ws.Range("E40").Value = "  -2.45%  "
$ws.Range("D41").Value = "0.9425"
$ws.Range("E41").Value = "  +3.79%  "
$ws.Range("E42").Value = "  +0.00%  "
$ws.Range("D43").Value = "1.989.84"
$ws.Range("E43").Value = "  -0.89%  "
$ws.Range("D44").Value = "100.53"
$ws.Range("E44").Value = "  -1.06%  "
$ws.Range("D45").Value = "65.40"
$ws.Range("E45").Value = "  -1.21%  "
$ws.Range("E46").Value = "  +2.01%  "
$ws.Range("D47").Value = "1.705"
$ws.Range("E47").Value = "  -2.31%  "
$ws.Range("D48").Value = "6.941"
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "8.892"
$ws.Range("E49").Value = "  -1.13%  "
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").Value = "0.1138"
$ws.Range("E50").Value = "  -2.87%  "
$ws.Range("E51").Value = "  -1.19%  "
